$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("assert_data")

$ws.Range("A3").Value = "item_header"
$ws.Range("B3").Value = "Canon EOS Rebel T7 DSLR Camera with 18-55mm IS Lens Kit"

$ws.Rows.Item(3).RowHeight = 18

$ws.Range("B3").Font.Color = 2893085

$ws.Range("B7").Select()
